# CS156 Assignment #3 Experiment Results — update the experiment-results
# table with the re-measured numbers, and let the trailing page-number
# field in the footer catch up to the now-longer document.

$d = $word.ActiveDocument

# The results table is the 2nd table in the document (1-based index in the
# Word object model): Metric | Type | Australia | Hardest4x4 | Easiest9x9 |
# Easy9x9 | Medium9x9 | Hard9x9 | Hardest9x9
$table = $d.Tables.Item(2)

function Set-CellText($tbl, $row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $cell.Range.Text = $text
}

# Row 2: Backtrack Function Call Count / With FC
Set-CellText $table 2 6 "136"
Set-CellText $table 2 7 "112"
Set-CellText $table 2 8 "83"
Set-CellText $table 2 9 "168"

# Row 3: Backtrack Function Call Count / No FC
Set-CellText $table 3 4 "77"
Set-CellText $table 3 5 "280"

# Row 4: Variable Value Setting Count / With FC
Set-CellText $table 4 6 "143"
Set-CellText $table 4 7 "121"
Set-CellText $table 4 8 "84"
Set-CellText $table 4 9 "176"

# Row 5: Variable Value Setting Count / No FC
Set-CellText $table 5 4 "278"
Set-CellText $table 5 5 "2014"

# Row 6: Algorithm Runtime (s) / With FC
Set-CellText $table 6 4 "0.016"
Set-CellText $table 6 5 "0.062"
Set-CellText $table 6 6 "0.094"
Set-CellText $table 6 7 "0.078"
Set-CellText $table 6 8 "0.063"
Set-CellText $table 6 9 "0.125s"

# Row 7: Algorithm Runtime (s) / No FC
Set-CellText $table 7 4 "0.016"
Set-CellText $table 7 5 "0.254"

# The document grew by a page once the table above was re-populated, so the
# cached result of the footer's PAGE field needs to move from 1 to 2.
$footer = $d.Sections.Item(1).Footers.Item(1)
$fields = $footer.Range.Fields
for ($i = 1; $i -le $fields.Count; $i++) {
    $fld = $fields.Item($i)
    if ($fld.Code.Text -match "PAGE") {
        $fld.Result.Text = "2"
    }
}
